$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet1" to "Mapping Samples"
$ws.Name = "Mapping Samples"

# Clear the empty separator row (row 5) contents and formatting entirely,
# without shifting the rows below it.
$ws.Range("A5:G5").Clear()

# Update the active selection on the sheet
$ws.Range("E35").Select()
